$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column before D ("Unit price"). This shifts the existing
#    Qty/Price/Notes columns (D,E,F) one place right (-> E,F,G) and Excel's
#    column-insert logic keeps formulas / hyperlinks / widths in sync.
# ---------------------------------------------------------------------------
$ws.Columns("D").Insert()

# New column header + widths.
$ws.Range("D1").Value = "Unit price"
$ws.Columns("D").ColumnWidth = 13.67
$ws.Columns("F").ColumnWidth = 20.83

# ---------------------------------------------------------------------------
# 2. Fill in the new "Unit price" values for the existing rows, and rebuild
#    the "Price" column as a real formula (Qty * Unit price).
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = 8.94
$ws.Range("D3").Value = 1.41
$ws.Range("D4").Value = 1.1399999999999999
$ws.Range("D5").Value = 6.58

$ws.Range("F2").Formula = "=E2*D2"
$ws.Range("F2").Style = "Normal"
$ws.Range("D2").Style = "Normal"

# F3:F7 share one formula group (built below, after rows 6-7 exist too).

# ---------------------------------------------------------------------------
# 3. New parts in the previously-empty rows 6 and 7.
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "POT 10K OHM LINEAR"
$ws.Range("B6").Value = "1993-1066-ND"
$ws.Range("C6").Value = "PT10MV10-103A2020-E-S"
$ws.Range("D6").Value = 0.49
$ws.Range("E6").Value = 4

$ws.Range("A7").Value = " Components"
$ws.Range("B7").Value = "CAP TRIMMER 2.5-22PF 250V TH"
$ws.Range("C7").Value = "BFC280800006"
$ws.Range("D7").Value = 2.89
$ws.Range("E7").Value = 2

# Shared "Qty * Unit price" formula for rows 3-7 (matches the ref="F3:F7" group).
$ws.Range("F3:F7").Formula = "=E3*D3"

# Apply the "Price" number/font style (style used by the original E2/F2 cell)
# to the whole Price column.
$ws.Range("F2:F7").Style = "Normal"
$ws.Range("D2").Font.Name = $ws.Range("D2").Font.Name
$ws.Range("A2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("D2").Value = 8.94
$ws.Range("F2").Copy()
$ws.Range("F2:F7").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. Styling that doesn't come from plain numbers/text:
#      - A6 / B7 use a new "black Arial 12" font.
#      - A7 / B6 / C6 / C7 reuse the existing Hyperlink look (same as G5).
# ---------------------------------------------------------------------------
$ws.Range("A3").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Font.Color = 0

$ws.Range("A6").Copy()
$ws.Range("B7").PasteSpecial(-4122)

$ws.Range("G5").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("A7").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5. Hyperlinks for the new parts (re-apply the Hyperlink look afterwards,
#    since adding a hyperlink resets formatting on the target cell).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B6"), "https://www.digikey.com/en/products/detail/bourns-inc/PT10MV10-103A2020-E-S/1750289")
$ws.Range("G5").Copy()
$ws.Range("B6").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("C6"), "https://www.digikey.com/en/products/detail/bourns-inc/PT10MV10-103A2020-E-S/1750289")
$ws.Range("G5").Copy()
$ws.Range("C6").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("C7"), "https://www.digikey.com/en/products/detail/vishay-bc-components/BFC280800006/285852")
$ws.Range("G5").Copy()
$ws.Range("C7").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("A7"), "https://www.digikey.com/en/products/detail/vishay-bc-components/BFC280800006/285852")
$ws.Range("G5").Copy()
$ws.Range("A7").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 6. Totals row: SUM now lives in column F (shifted from the old E40).
# ---------------------------------------------------------------------------
$ws.Range("F40").Formula = "=SUM(F2:F38)"

# ---------------------------------------------------------------------------
# 7. Selection / cursor position, like the saved workbook shows.
# ---------------------------------------------------------------------------
$ws.Range("E8").Select()
